# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking price strings are written with a leading apostrophe so Excel
# keeps storing them as text (matching the source data) instead of silently
# converting them to numbers and dropping formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.818.06'
$ws.Range("E2").Value = '  +1.62%  '
# Row 3
$ws.Range("D3").Value = '2.955.25'
$ws.Range("E3").Value = '  -0.90%  '
# Row 4
$ws.Range("E4").Value = '  -0.12%  '
# Row 5
$ws.Range("D5").Value = '''567.92'
$ws.Range("E5").Value = '  -2.18%  '
# Row 6
$ws.Range("D6").Value = '''160.62'
$ws.Range("E6").Value = '  +5.39%  '
# Row 7
$ws.Range("E7").Value = '  -0.10%  '
# Row 8
$ws.Range("E8").Value = '  +1.53%  '
# Row 9
$ws.Range("D9").Value = '2.951.79'
$ws.Range("E9").Value = '  -0.85%  '
# Row 10
$ws.Range("D10").Value = '''6.72'
$ws.Range("E10").Value = '  -3.24%  '
# Row 11
$ws.Range("E11").Value = '  -0.97%  '
# Row 12
$ws.Range("D12").Value = '''0.453'
$ws.Range("E12").Value = '  +1.66%  '
# Row 13
$ws.Range("E13").Value = '  +2.49%  '
# Row 14
$ws.Range("D14").Value = '''34.16'
$ws.Range("E14").Value = '  -0.03%  '
# Row 15
$ws.Range("E15").Value = '  -0.63%  '
# Row 16
$ws.Range("D16").Value = '65.912.35'
$ws.Range("E16").Value = '  +1.63%  '
# Row 17
$ws.Range("D17").Value = '3.444.94'
$ws.Range("E17").Value = '  -0.95%  '
# Row 18
$ws.Range("D18").Value = '''6.94'
$ws.Range("E18").Value = '  +0.74%  '
# Row 19
$ws.Range("D19").Value = '2.953.32'
$ws.Range("E19").Value = '  -1.04%  '
# Row 20
$ws.Range("D20").Value = '''444.69'
$ws.Range("E20").Value = '  -0.44%  '
# Row 21
$ws.Range("E21").Value = '  +0.80%  '
# Row 22
$ws.Range("E22").Value = '  -0.37%  '
# Row 23
$ws.Range("D23").Value = '''7.21'
$ws.Range("E23").Value = '  -1.07%  '
# Row 24
$ws.Range("D24").Value = '''82.18'
$ws.Range("E24").Value = '  +1.42%  '
# Row 25
$ws.Range("E25").Value = '  +0.58%  '
# Row 26
$ws.Range("D26").Value = '''12.16'
$ws.Range("E26").Value = '  -0.51%  '
# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  +0.02%  '
# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''9.97'
$ws.Range("E28").Value = '  -8.34%  '
# Row 29
$ws.Range("D29").Value = '''8.07'
$ws.Range("E29").Value = '  +4.85%  '
# Row 30
$ws.Range("E30").Value = '  -1.62%  '
# Row 31
$ws.Range("E31").Value = '  +0.02%  '
# Row 32
$ws.Range("D32").Value = '0.0₃0968'
$ws.Range("E32").Value = '  -9.28%  '
# Row 33
$ws.Range("D33").Value = '''27.15'
$ws.Range("E33").Value = '  +2.25%  '
# Row 34
$ws.Range("E34").Value = '  +0.51%  '
# Row 35
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  -0.16%  '
# Row 36
$ws.Range("D36").Value = '''0.977'
$ws.Range("E36").Value = '  -0.41%  '
# Row 37
$ws.Range("D37").Value = '''5.68'
$ws.Range("E37").Value = '  +0.75%  '
# Row 38
$ws.Range("D38").Value = '''49.21'
$ws.Range("E38").Value = '  +0.70%  '
# Row 39
$ws.Range("E39").Value = '  -6.24%  '
# Row 40
$ws.Range("E40").Value = '  +2.29%  '
# Row 41
$ws.Range("B41").Value = 'Arweave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D41").Value = '''43.39'
$ws.Range("E41").Value = '  -1.11%  '
# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.119'
$ws.Range("E42").Value = '  -1.22%  '
# Row 43
$ws.Range("E43").Value = '  -6.22%  '
# Row 44
$ws.Range("E44").Value = '  +0.01%  '
# Row 45
$ws.Range("D45").Value = '''380.95'
$ws.Range("E45").Value = '  -0.60%  '
# Row 46
$ws.Range("E46").Value = '  +1.56%  '
# Row 47
$ws.Range("D47").Value = '2.717.49'
$ws.Range("E47").Value = '  -2.45%  '
# Row 48
$ws.Range("D48").Value = '''130.39'
$ws.Range("E48").Value = '  -3.40%  '
# Row 49
$ws.Range("E49").Value = '  +0.05%  '
# Row 50
$ws.Range("E50").Value = '  +1.00%  '
# Row 51
$ws.Range("D51").Value = '''23.20'
$ws.Range("E51").Value = '  +1.44%  '
